# "Final Changes in Front-End, updated Back-End"
#
# Updates progress-tracking figures on the (only) worksheet of the
# To-do list workbook:
#   - Front-End Applikation / Unit-Tests (row 26): now fully done (100%,
#     status colour flips from "in progress" yellow to "done" green)
#   - Back-End Vorbereitung / Restschnittstellen (row 30): bumped to 100%
#   - Back-End Applikation / Code-Dokumentation (row 36): now fully done
#     (100%, status colour flips from yellow to green)
#   - Back-End Applikation / Unit-Tests (row 37): bumped to 50%
# The dependent percentage-summary formulas in column C recalculate
# automatically.
# Also moves the on-screen selection/scroll position to where the author
# left off working (around row 13, cell H28 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Front-End Applikation section ---------------------------------
# Unit-Tests: 25 -> 100, status fill yellow -> green
$ws.Range("B26").Interior.Color = 5296274   # RGB(146,208,80) "done" green
$ws.Range("C26").Value = 100

# --- Back-End Vorbereitung section ----------------------------------
# Restschnittstellen: 99 -> 100
$ws.Range("C30").Value = 100

# --- Back-End Applikation section -----------------------------------
# Code-Dokumentation: 20 -> 100, status fill yellow -> green
$ws.Range("B36").Interior.Color = 5296274   # RGB(146,208,80) "done" green
$ws.Range("C36").Value = 100

# Unit-Tests: 20 -> 50
$ws.Range("C37").Value = 50

# --- View state: scroll to row 13 and select H28 --------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("H28").Select()
